# Add a new "Other Art" worksheet before the first sheet ("Animations")
$wb = $excel.ActiveWorkbook
$other = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$other.Name = "Other Art"

$other.Range("A1").Value = "Selector"
$other.Range("A2").Value = "Fog"
$other.Range("A3").Value = "Select Creature Screen"
$other.Range("A4").Value = "Main Menu Screen"
$other.Range("A5").Value = "Options Screen"
$other.Range("A6").Value = "Pause Menu"
$other.Range("A7").Value = "Level Backgrounds (5)"
$other.Columns.Item(1).ColumnWidth = 12.7

# Remove the now-redundant "Level Backgrounds" sheet -- consolidated into
# the "Other Art" list above ("Level Backgrounds (5)")
$wb.Worksheets.Item("Level Backgrounds").Delete() | Out-Null

# A couple of new Animations entries now have an associated ($) cost, so
# tag those cells with a currency number format
$animations = $wb.Worksheets.Item("Animations")
$animations.Range("H3").NumberFormat = "$#,##0_);[Red]($#,##0)"
$animations.Range("H4").NumberFormat = "$#,##0_);[Red]($#,##0)"
$animations.Range("G16").Select() | Out-Null

# Update the Sound sheet's selection and make "Other Art" the active tab
$sound = $wb.Worksheets.Item("Sound")
$sound.Range("E14").Select() | Out-Null

$other.Activate() | Out-Null
$other.Range("A8").Select() | Out-Null
